$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F7").Value = 78
$ws.Range("G7").Value = 3652.74
$ws.Range("F8").Value = 62
$ws.Range("G8").Value = 14432.36
$ws.Range("B9").Value = 18491.79
$ws.Range("F23").Value = 14
$ws.Range("G23").Value = 573.86
$ws.Range("F24").Value = 26
$ws.Range("G24").Value = 2663.7
$ws.Range("F27").Value = 59
$ws.Range("G27").Value = 2115.74
$ws.Range("F29").Value = 75
$ws.Range("G29").Value = 3842.25
$ws.Range("B34").Value = 60576.4
$ws.Range("F36").Value = 93
$ws.Range("G36").Value = 18299.61
$ws.Range("F41").Value = 217
$ws.Range("G41").Value = 41857.13
$ws.Range("F42").Value = 68
$ws.Range("G42").Value = 2864.16
$ws.Range("F46").Value = 59
$ws.Range("G46").Value = 2193.62
$ws.Range("F53").Value = 35
$ws.Range("G53").Value = 574.35
$ws.Range("F55").Value = 130
$ws.Range("G55").Value = 7248.8
$ws.Range("F56").Value = 41
$ws.Range("G56").Value = 915.12
$ws.Range("F58").Value = 80
$ws.Range("G58").Value = 6234.4
$ws.Range("F61").Value = 237
$ws.Range("G61").Value = 61793.01
$ws.Range("F65").Value = 14
$ws.Range("G65").Value = 468.44
$ws.Range("B66").Value = 213448.48
$ws.Range("F86").Value = 2
$ws.Range("G86").Value = 14008.92
$ws.Range("B87").Value = 14008.92
$ws.Range("F143").Value = 44
$ws.Range("G143").Value = 8135.16
$ws.Range("B147").Value = 22530.65
$ws.Range("F167").Value = 0
$ws.Range("G167").Value = 0
$ws.Range("B170").Value = 11366.63
$ws.Range("F186").Value = 24
$ws.Range("G186").Value = 1038.72
$ws.Range("B193").Value = 68253.87
$ws.Range("F215").Value = 176
$ws.Range("G215").Value = 19763.04
$ws.Range("B218").Value = 83380.47
$ws.Range("F222").Value = 1005
$ws.Range("G222").Value = 18592.5
$ws.Range("F227").Value = 51
$ws.Range("G227").Value = 5844.6
$ws.Range("B229").Value = 31561.21
$ws.Range("F232").Value = 27
$ws.Range("G232").Value = 3095.01
$ws.Range("B240").Value = 14879.7
$ws.Range("F255").Value = 13
$ws.Range("G255").Value = 4105.4
$ws.Range("F268").Value = 12
$ws.Range("G268").Value = 1526.52
$ws.Range("F277").Value = 12
$ws.Range("G277").Value = 604.92
$ws.Range("F278").Value = 42
$ws.Range("G278").Value = 5692.68
$ws.Range("F287").Value = 62
$ws.Range("G287").Value = 3393.88
$ws.Range("B295").Value = 126770.45
$ws.Range("B297").Value = 61610
$ws.Range("E297").Value = 122.71
$ws.Range("F297").Value = -58
$ws.Range("G297").Value = -5957.18
$ws.Range("B298").Value = 63565
$ws.Range("E298").Value = 109.19
$ws.Range("F298").Value = 60
$ws.Range("G298").Value = 6162.6
$ws.Range("B306").Value = 63531
$ws.Range("E306").Value = 152.53
$ws.Range("F306").Value = 29
$ws.Range("G306").Value = 4160.92
$ws.Range("B307").Value = 57802
$ws.Range("E307").Value = 162.71
$ws.Range("F307").Value = -79
$ws.Range("G307").Value = -11334.92
$ws.Range("F324").Value = 57
$ws.Range("G324").Value = 9765.809999999999
$ws.Range("B328").Value = -1437.97
$ws.Range("F366").Value = 65
$ws.Range("G366").Value = 3596.45
$ws.Range("F371").Value = 69
$ws.Range("G371").Value = 10368.63
$ws.Range("B372").Value = 65167.27
$ws.Range("F387").Value = 441
$ws.Range("G387").Value = 42600.6
$ws.Range("B389").Value = 59542.06
$ws.Range("F396").Value = 136
$ws.Range("G396").Value = 3465.28
$ws.Range("F402").Value = 53
$ws.Range("G402").Value = 1818.43
$ws.Range("F403").Value = 78
$ws.Range("G403").Value = 3162.12
$ws.Range("F415").Value = 60
$ws.Range("G415").Value = 3282
$ws.Range("B417").Value = 175096.98
$ws.Range("F433").Value = 142
$ws.Range("G433").Value = 1368.88
$ws.Range("B438").Value = 26523.44
$ws.Range("F450").Value = 12
$ws.Range("G450").Value = 2666.76
$ws.Range("F454").Value = 77
$ws.Range("G454").Value = 21787.92
$ws.Range("B458").Value = 101034.53
$ws.Range("B479").Value = 64810
$ws.Range("E479").Value = 291.22
$ws.Range("F479").Value = 0
$ws.Range("G479").Value = 0
$ws.Range("B480").Value = 53319
$ws.Range("E480").Value = 310.64
$ws.Range("F480").Value = -6
$ws.Range("G480").Value = -1643.52
$ws.Range("F513").Value = 223
$ws.Range("G513").Value = 11031.81
$ws.Range("B525").Value = 130864.33
$ws.Range("F533").Value = 3
$ws.Range("G533").Value = 90.23999999999999
$ws.Range("B535").Value = 26146.12
$ws.Range("F544").Value = 45
$ws.Range("G544").Value = 2785.5
$ws.Range("F545").Value = 2
$ws.Range("G545").Value = 5771.4
$ws.Range("F551").Value = 12
$ws.Range("G551").Value = 8942.280000000001
$ws.Range("F552").Value = 0
$ws.Range("G552").Value = 0
$ws.Range("B556").Value = 54417.98
$ws.Range("F558").Value = 216
$ws.Range("G558").Value = 26319.6
$ws.Range("B561").Value = 31131.32
$ws.Range("F569").Value = 4
$ws.Range("G569").Value = 2338.88
$ws.Range("B573").Value = 28771.2
$ws.Range("F605").Value = 197
$ws.Range("G605").Value = 26220.7
$ws.Range("B607").Value = 26760.74
$ws.Range("F609").Value = 26
$ws.Range("G609").Value = 2829.06
$ws.Range("F612").Value = 237
$ws.Range("G612").Value = 35647.17
$ws.Range("F623").Value = 82
$ws.Range("G623").Value = 42192.28
$ws.Range("F626").Value = 13
$ws.Range("G626").Value = 613.73
$ws.Range("B628").Value = 217637.69
$ws.Range("F659").Value = 43
$ws.Range("G659").Value = 2302.22
$ws.Range("F662").Value = 49
$ws.Range("G662").Value = 3935.19
$ws.Range("B668").Value = 13360.89
$ws.Range("F674").Value = 916
$ws.Range("G674").Value = 149408.76
$ws.Range("B680").Value = 150421.31
$ws.Range("B718").Value = 2853457.61
$ws.Range("B719").Value = 2853457.61
